$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header label (K1) describing the weekly progress grid
$ws.Range("K1").Value = "Semanas com o progresso de cada passo da atividade"

# Week numbers 1-9 across K4:S4 (one column per week)
$weekRange = $ws.Range("K4:S4")
for ($i = 0; $i -lt 9; $i++) {
    $weekRange.Cells.Item(1, $i + 1).Value = $i + 1
}

# Restore the saved selection/active cell
$ws.Range("O8").Select()
